# Updating filtered feeds from workflow
# Appends two new feed rows (GenomeWeb + 360Dx mirror) about Pillar Biosciences'
# CMS coverage for its pan-cancer CDx assay, tagged with the "CDx" keyword.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$genomewebLink = "https://www.genomeweb.com/cancer/pillar-biosciences-receives-cms-coverage-pan-cancer-cdx-assay"
$dxLink        = "https://www.360dx.com/cancer/pillar-biosciences-receives-cms-coverage-pan-cancer-cdx-assay"
$keyword       = "CDx"
$title         = "Pillar Biosciences Receives CMS Coverage for Pan-Cancer CDx Assay"

# Row 40: GenomeWeb article
$ws.Range("A40").Value = $genomewebLink
$ws.Range("B40").Value = $keyword
$ws.Range("C40").Value = $title
$ws.Hyperlinks.Add($ws.Range("A40"), $genomewebLink)
$ws.Range("A40").Style = $ws.Range("A39").Style

# Row 41: 360Dx mirror article
$ws.Range("A41").Value = $dxLink
$ws.Range("B41").Value = $keyword
$ws.Range("C41").Value = $title
$ws.Hyperlinks.Add($ws.Range("A41"), $dxLink)
$ws.Range("A41").Style = $ws.Range("A39").Style
